# Update power-up credits and enhance team data structure in Excel export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up a few "Users" entries (sharedStrings content updates) ---
$ws.Range("C2").Value  = "Ankush Gautam, Vaibhav Srivastva, utkal, Aishlee Joshi"
$ws.Range("C25").Value = "Sameer Verma, Aditya, Palak, Vaibhav Gupta"
$ws.Range("C26").Value = "Aman, Angad, harshil, Kashish"
$ws.Range("C28").Value = "Aakarsh, Aryan2, Anirudh, Pranav7"

# --- 2. Add new "Score" column (J) ---
$ws.Range("J1").Value = "Score"

$scores = @(2835,1418,3875,2930,4347,7182,5103,5387,5935,2268,3591,5954,2457,2740,6898,1985,0,1985,4309,3875,2268,6520,1059,3780,3100,662,2169,2741)

for ($i = 0; $i -lt $scores.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $scores[$i]
}

# Match the column width convention used by the other columns (XML width "10")
$ws.Columns.Item(10).ColumnWidth = 9.17

Write-Host "Applied team data updates and Score column."
